# Change "Point" column (A) to "Time": header text + numeric time values
# instead of P1..P20 string labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("A1").Value = "Time"

# New numeric values for A2:A21 (replacing P1..P20 text labels)
$times = @(1, 3, 4, 6, 9, 12, 14, 15, 17, 18, 20, 21, 25, 27, 28, 30, 31, 32, 38, 40)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $times[$i]
}

# Update selection / scroll position to match the new view
$ws.Range("I13").Select()
